$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "备注"
$ws.Range("C2").Value = "1组"

$ws.Range("C2").Select()
